$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: D0.5313897347289666 -> D0.9125181743180528
$ws.Range("A2").Value = "D0.9125181743180528"

# A3: D0.6763565034871812 -> D0.5477877949969697
$ws.Range("A3").Value = "D0.5477877949969697"

# A4: D0.742386377676004 -> D0.6037163037954872
$ws.Range("A4").Value = "D0.6037163037954872"

# D2, D3, D4 all share the same string: Sat, 10 Dec 2022 21:22:41 +0530 -> Fri, 23 Dec 2022 23:19:07 -0800
$ws.Range("D2").Value = "Fri, 23 Dec 2022 23:19:07 -0800"
$ws.Range("D3").Value = "Fri, 23 Dec 2022 23:19:07 -0800"
$ws.Range("D4").Value = "Fri, 23 Dec 2022 23:19:07 -0800"
